$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns I and J, same style as the existing H1 header ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# --- Data rows 2-69: new numeric columns I (I0) and J (IF) ---
$iVals = @(7,9,7,5,6,6,9,8,6,7,7,8,8,9,6,7,7,6,5,7,7,7,7,9,6,9,8,8,7,7,8,7,6,6,9,9,7,8,7,6,7,8,9,9,6,9,8,8,8,8,7,6,8,9,9,8,8,7,10,8,7,8,8,8,8,7,5,8)
$jVals = @(7,9,7,6,6,7,9,8,6,7,7,8,9,9,6,7,7,6,6,7,8,7,7,9,7,9,9,9,7,7,8,7,7,7,9,9,7,8,7,7,8,8,9,9,7,9,8,8,9,8,7,7,8,9,9,8,8,8,10,8,7,8,8,8,8,7,6,8)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value2 = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value2 = $jVals[$idx]
}

Write-Output ("Dimension: " + $ws.UsedRange.Address())
